$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Matches_SOG: append 5 new match rows (465-469)
# ---------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @("897763", "2025-11-12T17:00:00", "Салават Юлаев", "СКА", 30, 29, "khl_text"),
    @("897766", "2025-11-12T17:00:00", "Металлург Мг", "Лада", 42, 14, "khl_text"),
    @("897764", "2025-11-12T19:30:00", "Торпедо", "ХК Сочи", 43, 24, "khl_text"),
    @("897765", "2025-11-12T19:30:00", "Динамо Мн", "Сибирь", 55, 20, "khl_text"),
    @("897762", "2025-11-12T19:30:00", "Спартак", "Локомотив", 27, 32, "khl_text")
)

$r = 465
foreach ($row in $newMatches) {
    # Column A (uid) looks numeric ("897763") - force text so it matches the
    # existing uid column storage (text, not a number).
    $wsMatches.Range("A$r").NumberFormat = "@"
    $wsMatches.Cells.Item($r, 1).Value = $row[0]
    $wsMatches.Cells.Item($r, 2).Value = $row[1]
    $wsMatches.Cells.Item($r, 3).Value = $row[2]
    $wsMatches.Cells.Item($r, 4).Value = $row[3]
    $wsMatches.Cells.Item($r, 5).Value = $row[4]
    $wsMatches.Cells.Item($r, 6).Value = $row[5]
    $wsMatches.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------
# 2) Shots_HA: refresh as_of_utc + recomputed home/away shot stats
# ---------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$haData = @{
    2 = @(22, 18, 730, 616, 33.2, 28, 642, 550, 35.7, 30.6)  # Авангард
    3 = @(18, 27, 521, 558, 28.9, 31, 752, 831, 27.9, 30.8)  # Автомобилист
    4 = @(18, 20, 686, 483, 38.1, 26.8, 638, 561, 31.9, 28.1)  # Адмирал
    5 = @(23, 21, 779, 587, 33.9, 25.5, 711, 627, 33.9, 29.9)  # Ак Барс
    6 = @(21, 21, 641, 746, 30.5, 35.5, 592, 769, 28.2, 36.6)  # Амур
    7 = @(30, 15, 957, 953, 31.9, 31.8, 415, 515, 27.7, 34.3)  # Барыс
    8 = @(18, 21, 598, 491, 33.2, 27.3, 596, 673, 28.4, 32)  # Динамо М
    9 = @(25, 18, 920, 677, 36.8, 27.1, 668, 484, 37.1, 26.9)  # Динамо Мн
    10 = @(20, 21, 571, 698, 28.6, 34.9, 573, 772, 27.3, 36.8)  # Драконы
    11 = @(22, 21, 598, 776, 27.2, 35.3, 531, 797, 25.3, 38)  # Лада
    12 = @(18, 27, 555, 495, 30.8, 27.5, 839, 666, 31.1, 24.7)  # Локомотив
    13 = @(26, 17, 926, 660, 35.6, 25.4, 488, 460, 28.7, 27.1)  # Металлург Мг
    14 = @(25, 19, 781, 858, 31.2, 34.3, 517, 722, 27.2, 38)  # Нефтехимик
    15 = @(25, 18, 830, 834, 33.2, 33.4, 554, 590, 30.8, 32.8)  # СКА
    16 = @(17, 26, 465, 469, 27.4, 27.6, 730, 758, 28.1, 29.2)  # Салават Юлаев
    17 = @(16, 25, 469, 363, 29.3, 22.7, 829, 665, 33.2, 26.6)  # Северсталь
    18 = @(19, 24, 511, 733, 26.9, 38.6, 657, 782, 27.4, 32.6)  # Сибирь
    19 = @(26, 15, 927, 721, 35.7, 27.7, 518, 534, 34.5, 35.6)  # Спартак
    20 = @(22, 27, 717, 650, 32.6, 29.5, 931, 866, 34.5, 32.1)  # Торпедо
    21 = @(18, 26, 599, 541, 33.3, 30.1, 888, 847, 34.2, 32.6)  # Трактор
    22 = @(19, 20, 560, 622, 29.5, 32.7, 516, 745, 25.8, 37.2)  # ХК Сочи
    23 = @(20, 21, 470, 577, 23.5, 28.9, 523, 597, 24.9, 28.4)  # ЦСКА
}

foreach ($row in 2..23) {
    $wsHA.Cells.Item($row, 4).Value = "2025-11-12T19:30:00Z"
    $vals = $haData[$row]
    $wsHA.Cells.Item($row, 5).Value = $vals[0]
    $wsHA.Cells.Item($row, 6).Value = $vals[1]
    $wsHA.Cells.Item($row, 7).Value = $vals[2]
    $wsHA.Cells.Item($row, 8).Value = $vals[3]
    $wsHA.Cells.Item($row, 9).Value = $vals[4]
    $wsHA.Cells.Item($row, 10).Value = $vals[5]
    $wsHA.Cells.Item($row, 11).Value = $vals[6]
    $wsHA.Cells.Item($row, 12).Value = $vals[7]
    $wsHA.Cells.Item($row, 13).Value = $vals[8]
    $wsHA.Cells.Item($row, 14).Value = $vals[9]
}

# ---------------------------------------------------------------
# 3) Shots_Summary: refresh as_of_utc + recomputed SOG totals
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

$summaryData = @{
    2 = @(40, 1372, 1166, 34.3, 29.1)  # Авангард
    3 = @(45, 1273, 1389, 28.3, 30.9)  # Автомобилист
    4 = @(38, 1324, 1044, 34.8, 27.5)  # Адмирал
    5 = @(44, 1490, 1214, 33.9, 27.6)  # Ак Барс
    6 = @(42, 1233, 1515, 29.4, 36.1)  # Амур
    7 = @(45, 1372, 1468, 30.5, 32.6)  # Барыс
    8 = @(39, 1194, 1164, 30.6, 29.8)  # Динамо М
    9 = @(43, 1588, 1161, 36.9, 27)  # Динамо Мн
    10 = @(41, 1144, 1470, 27.9, 35.9)  # Драконы
    11 = @(43, 1129, 1573, 26.3, 36.6)  # Лада
    12 = @(45, 1394, 1161, 31, 25.8)  # Локомотив
    13 = @(43, 1414, 1120, 32.9, 26)  # Металлург Мг
    14 = @(44, 1298, 1580, 29.5, 35.9)  # Нефтехимик
    15 = @(43, 1384, 1424, 32.2, 33.1)  # СКА
    16 = @(43, 1195, 1227, 27.8, 28.5)  # Салават Юлаев
    17 = @(41, 1298, 1028, 31.7, 25.1)  # Северсталь
    18 = @(43, 1168, 1515, 27.2, 35.2)  # Сибирь
    19 = @(41, 1445, 1255, 35.2, 30.6)  # Спартак
    20 = @(49, 1648, 1516, 33.6, 30.9)  # Торпедо
    21 = @(44, 1487, 1388, 33.8, 31.5)  # Трактор
    22 = @(39, 1076, 1367, 27.6, 35.1)  # ХК Сочи
    23 = @(41, 993, 1174, 24.2, 28.6)  # ЦСКА
}

foreach ($row in 2..23) {
    $wsSummary.Cells.Item($row, 4).Value = "2025-11-12T19:30:00Z"
    $vals = $summaryData[$row]
    $wsSummary.Cells.Item($row, 5).Value = $vals[0]
    $wsSummary.Cells.Item($row, 6).Value = $vals[1]
    $wsSummary.Cells.Item($row, 7).Value = $vals[2]
    $wsSummary.Cells.Item($row, 8).Value = $vals[3]
    $wsSummary.Cells.Item($row, 9).Value = $vals[4]
}

# ---------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-12T19:30:00Z"
$wsMeta.Range("D2").Value = 61
